# Update 28-Feb-2021, end of day.
# Clears out the detailed petty-cash transactions for the old date range
# (rows 4-34) and the expense-type formula in D3, resets the opening
# balance (E2) and the date of the first remaining entry (A3), then
# leaves the selection on D4 ready for the next day's entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New opening balance carried forward.
$ws.Range("E2").Value = 687025

# First row keeps its date (updated) and description, but its Debit
# formula is removed.
$ws.Range("A3").Value = 44256
$ws.Range("D3").Clear()

# All the other transaction rows (date/description/debit/credit) for
# this block are cleared out entirely.
$ws.Range("A4:D34").Clear()

# Leave the selection where the next entry would be typed.
$ws.Range("D4").Select()
